$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fix the misspelled "Accpted" status entries to "Accepted" (but not the G8 cell
#    which already holds a different, unrelated value).
$statusRange = $ws.Range("G3:G10")
foreach ($cell in $statusRange.Cells) {
    if ($cell.Value() -eq "Accpted") {
        $cell.Value = "Accepted"
    }
}

# 2. Update the Project Schedule review comment (cell E4) with the revised wording.
$newComment = "1- Review the start dates and finish dates with the  team members.`n2- Rename the section of 'work' to 'work hours'.`n3- The Completed percentage section calculations is incorrect.`n4- Move the question marks displayed after the number of days from the 'duration' section.`n5- The extension of 'TAWA_SIQ 'document is '.xlsx 'not '.docx'.`n6- Add the' TAWA_SystemRequirements.xlsx' document and 'TAWA_CustomerRequirements.xlsx' document to the documents section.`n7- The PDF version contents of the 'project schedule' document are overlapping .`n"
$ws.Range("E4").Value = $newComment

# 3. The comment in E4 is now longer, so grow row 4 to fit the extra lines.
$ws.Rows.Item(4).RowHeight = 250.2

# 4. Update the sheet view so it is scrolled down a bit with E4 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E4").Select()
